$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking strings in column D stay as text (preserve exact formatting,
# trailing zeros, multi-dot grouping, etc.) by forcing Text number format first.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Update price (D) and 1h volume change (E) values for rows 2-47 ---
$ws.Range("D2").Value = "26.137.18"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.656.51"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").Value = "219.00"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "0.5247"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("D8").Value = "0.2622"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").Value = "0.06296"
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").Value = "20.60"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "0.07795"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "4.494"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").Value = "1.662.06"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "1.884.26"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "0.5552"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "0.0₅7997"
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("D18").Value = "26.158.16"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "4.643"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").Value = "195.47"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").Value = "10.11"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").Value = "5.961"
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("D25").Value = "146.63"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").Value = "0.1205"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").Value = "7.165"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").Value = "1.495"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").Value = "0.05711"
$ws.Range("E30").Value = "  -2.97%  "
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "3.488"
$ws.Range("E32").Value = "  -2.75%  "
$ws.Range("D33").Value = "3.347"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").Value = "1.587"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D36").Value = "0.9527"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("D37").Value = "2.419"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "0.5693"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").Value = "5.953"
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("D41").Value = "1.059.80"
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("D42").Value = "0.8452"
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("D44").Value = "103.50"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("D45").Value = "1.795.14"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").Value = "57.87"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").Value = "0.05413"
$ws.Range("E47").Value = "  +4.90%  "

# --- Rows 48-50 reordered: BabyDogeCoin rises above Frax and Mantle ---
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₈105"
$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").Value = "1.008"
$ws.Range("E49").Value = "  -0.35%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.4400"
$ws.Range("E50").Value = "  +0.39%  "

# --- Row 51: EnergySwap price/volume update ---
$ws.Range("D51").Value = "8.031"
$ws.Range("E51").Value = "  +0.12%  "
